$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from column J (row 3-6) into new column K so the new
# cells inherit the same number formats / borders / fonts as the rest
# of the table, then fill in the 2023 values.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1215.9000000000001
$ws.Range("K5").Value = 1165.0999999999999
$ws.Range("K6").Value = 1300.5999999999999
